$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.844.01'
$ws.Range('E2').Value = '  +1.19%  '
$ws.Range('D3').Value = '1.618.04'
$ws.Range('D5').Value = "'213.29"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').Value = "'0.520"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D7').Value = "'0.992"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.89%  '
$ws.Range('D8').Value = "'29.26"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.22%  '
$ws.Range('E9').Value = '  +3.28%  '
$ws.Range('D10').Value = "'0.0606"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.65%  '
$ws.Range('D11').Value = "'0.0910"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').Value = '1.850.52'
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('D13').Value = '1.620.31'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = "'0.565"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.58%  '
$ws.Range('D15').Value = "'3.89"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.80%  '
$ws.Range('D16').Value = '29.865.61'
$ws.Range('E16').Value = '  +1.22%  '
$ws.Range('D17').Value = "'8.82"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +15.42%  '
$ws.Range('D18').Value = "'64.36"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.45%  '
$ws.Range('D19').Value = "'241.03"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = '0.0₃0706'
$ws.Range('E20').Value = '  +2.46%  '
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('D22').Value = "'4.09"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.41%  '
$ws.Range('D23').Value = "'9.58"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.09%  '
$ws.Range('D24').Value = "'2.11"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.20%  '
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').Value = "'6.57"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.93%  '
$ws.Range('E29').Value = '  -0.79%  '
$ws.Range('E30').Value = '  +3.03%  '
$ws.Range('D31').Value = "'1.12"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.66%  '
$ws.Range('D32').Value = "'3.34"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.56%  '
$ws.Range('E33').Value = '  +3.48%  '
$ws.Range('D34').Value = '1.415.38'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = "'1.63"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.10%  '
$ws.Range('D36').Value = "'1.03"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').Value = "'2.89"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.13%  '
$ws.Range('D38').Value = "'2.28"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('D40').Value = "'0.555"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.26%  '
$ws.Range('D41').Value = "'0.0502"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.99%  '
$ws.Range('D44').Value = "'53.56"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.10%  '
$ws.Range('D45').Value = "'69.17"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.12%  '
$ws.Range('E46').Value = '  +19.05%  '
$ws.Range('D47').Value = "'0.992"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('D48').Value = "'5.43"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.89%  '
$ws.Range('D49').Value = '1.759.32'
$ws.Range('E49').Value = '  +0.54%  '
$ws.Range('D50').Value = "'88.24"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.89%  '
$ws.Range('D51').Value = "'0.0533"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.08%  '

# Rows 42 and 43 swap content (ARBITRUM / RenderToken) with new values
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = "'0.825"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.47%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = "'1.97"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.22%  '
